$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("w01_100.5_optimization_results")

$ws.Range("B2").Value = 636
$ws.Range("C2").Value = 1.030354795083571
$ws.Range("E2").Value = -1
$ws.Range("G2").Value = 0.001

$ws.Range("B3").Value = 624
$ws.Range("D3").Value = 0.8722874034249394
$ws.Range("E3").Value = -0.1988943718412705
$ws.Range("F3").Value = 2.497255835710833
$ws.Range("G3").Value = 0.931135295829163

$ws.Range("B4").Value = 637
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = -1
$ws.Range("F4").Value = 3
$ws.Range("G4").Value = 0.001

$ws.Range("B5").Value = 638
$ws.Range("F5").Value = 2.380531664662494
$ws.Range("G5").Value = 0.001

$ws.Range("B6").Value = 638
$ws.Range("C6").Value = 0.6873389414318158
$ws.Range("D6").Value = 0.001
$ws.Range("E6").Value = -1
$ws.Range("G6").Value = 0.001

$ws.Range("B7").Value = 625
$ws.Range("C7").Value = 1.544792614611256
$ws.Range("D7").Value = 0.001
$ws.Range("F7").Value = 1.001
$ws.Range("G7").Value = 1

$ws.Range("B8").Value = 627
$ws.Range("F8").Value = 1.001
$ws.Range("G8").Value = 0.001

$ws.Range("B9").Value = 618
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = -1
$ws.Range("F9").Value = 3
$ws.Range("G9").Value = 0.001

$ws.Range("B10").Value = 627
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = -0.001
$ws.Range("F10").Value = 3

$ws.Range("B11").Value = 633
$ws.Range("C11").Value = 1.566938059915472
$ws.Range("D11").Value = 0.001
$ws.Range("E11").Value = -1
$ws.Range("F11").Value = 1.001
$ws.Range("G11").Value = 0.001
